# feat(objective): add secret objective to finish the camp
#
# Inserts a new localization row (key / en / fr) right after the
# "OBJECTIVE_THE_END" row (row 29) and before the "PAUSE_TITLE" row
# (old row 30), pushing every following row down by one. The new row
# carries the OBJECTIVE_GAME_FINISHED secret-objective strings. Column B
# ("en") is also widened to fit the new, longer English copy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row 30 (PAUSE_TITLE and everything below shifts to 31+).
$ws.Rows(30).Insert() | Out-Null

# Fill the new row with the secret "camp finished" objective strings.
$ws.Range("A30").Value = "OBJECTIVE_GAME_FINISHED"
$ws.Range("B30").Value = "You stayed strong! Congratulations, you finished the camp!"
$ws.Range("C30").Value = "Vous n’avez rien lâché ! Félicitions, vous avez fini votre campement !"

# Widen column B ("en") to comfortably fit the longer English text.
$ws.Columns(2).ColumnWidth = 49.67

# Match the author's final selection/cursor position.
$ws.Range("B33").Select() | Out-Null
